$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-40 with the corrected daily summary data.
# (Rows 41-42 are left untouched, matching the source diff.)

$ws.Cells.Item(2, 1).Value = 45908
$ws.Cells.Item(2, 2).Value = 'ARA3A'
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 4).Value = 57
$ws.Cells.Item(2, 5).Value = 'Benito A, Daniel GS, Evaristo A, Leobardo RL, Moises P, Omar S'

$ws.Cells.Item(3, 1).Value = 45908
$ws.Cells.Item(3, 2).Value = 'Founders 2'
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 28.5
$ws.Cells.Item(3, 5).Value = 'Honorio G, Rigoberto Al-B, Rogelio M'

$ws.Cells.Item(4, 1).Value = 45908
$ws.Cells.Item(4, 2).Value = 'GU Henle'
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = 38
$ws.Cells.Item(4, 5).Value = 'Fernando V, Laurentino, Noe VL, Oscar VS'

$ws.Cells.Item(5, 1).Value = 45908
$ws.Cells.Item(5, 2).Value = 'HanoverSpring'
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 57
$ws.Cells.Item(5, 5).Value = 'Antoine F, Carlos Al-V, Jose P, Rata F'

$ws.Cells.Item(6, 1).Value = 45908
$ws.Cells.Item(6, 2).Value = 'Rowan'
$ws.Cells.Item(6, 3).Value = 9
$ws.Cells.Item(6, 4).Value = 85.5
$ws.Cells.Item(6, 5).Value = 'Adalberto T, Daniel LG, Elvis T, Gaudencio B, Henry G, Julio M, Luis Enrique R, Luis Martin R, Trinidad T'

$ws.Cells.Item(7, 1).Value = 45908
$ws.Cells.Item(7, 2).Value = 'Tidal Basin'
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 4).Value = 57
$ws.Cells.Item(7, 5).Value = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR'

$ws.Cells.Item(8, 1).Value = 45908
$ws.Cells.Item(8, 2).Value = 'Wardman'
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 133
$ws.Cells.Item(8, 5).Value = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Jose Luis H, Juan G, Miguel A, Misael M, Pablo G, William A'

$ws.Cells.Item(9, 1).Value = 45909
$ws.Cells.Item(9, 2).Value = 'ARA3A'
$ws.Cells.Item(9, 3).Value = 6
$ws.Cells.Item(9, 4).Value = 57
$ws.Cells.Item(9, 5).Value = 'Benito A, Daniel GS, Eric M R, Evaristo A, Leobardo RL, Omar S'

$ws.Cells.Item(10, 1).Value = 45909
$ws.Cells.Item(10, 2).Value = 'Founders 2'
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 27
$ws.Cells.Item(10, 5).Value = 'Honorio G, Rigoberto Al-B, Rogelio M'

$ws.Cells.Item(11, 1).Value = 45909
$ws.Cells.Item(11, 2).Value = 'GU Henle'
$ws.Cells.Item(11, 3).Value = 4
$ws.Cells.Item(11, 4).Value = 38
$ws.Cells.Item(11, 5).Value = 'Fernando V, Laurentino, Noe VL, Oscar VS'

$ws.Cells.Item(12, 1).Value = 45909
$ws.Cells.Item(12, 2).Value = 'HanoverSpring'
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(12, 4).Value = 57
$ws.Cells.Item(12, 5).Value = 'Antoine F, Carlos Al-V, Gaudencio B, Jose P, Moises P, Rata F'

$ws.Cells.Item(13, 1).Value = 45909
$ws.Cells.Item(13, 2).Value = 'Rowan'
$ws.Cells.Item(13, 3).Value = 8
$ws.Cells.Item(13, 4).Value = 76
$ws.Cells.Item(13, 5).Value = 'Adalberto T, Daniel LG, Elvis T, Henry G, Julio M, Luis Enrique R, Luis Martin R, Trinidad T'

$ws.Cells.Item(14, 1).Value = 45909
$ws.Cells.Item(14, 2).Value = 'Tidal Basin'
$ws.Cells.Item(14, 3).Value = 6
$ws.Cells.Item(14, 4).Value = 57
$ws.Cells.Item(14, 5).Value = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR'

$ws.Cells.Item(15, 1).Value = 45909
$ws.Cells.Item(15, 2).Value = 'Wardman'
$ws.Cells.Item(15, 3).Value = 14
$ws.Cells.Item(15, 4).Value = 133
$ws.Cells.Item(15, 5).Value = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Jose Luis H, Juan G, Miguel A, Misael M, Pablo G, William A'

$ws.Cells.Item(16, 1).Value = 45910
$ws.Cells.Item(16, 2).Value = 'ARA3A     Moorefield'
$ws.Cells.Item(16, 3).Value = 7
$ws.Cells.Item(16, 4).Value = 67
$ws.Cells.Item(16, 5).Value = 'Benito A, Daniel GS, Eric M R, Evaristo A, Henry G, Leobardo RL, Omar S'

$ws.Cells.Item(17, 1).Value = 45910
$ws.Cells.Item(17, 2).Value = 'Canvas'
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 21
$ws.Cells.Item(17, 5).Value = 'Jose Carlos G, Noe VL'

$ws.Cells.Item(18, 1).Value = 45910
$ws.Cells.Item(18, 2).Value = 'Founders 2'
$ws.Cells.Item(18, 3).Value = 4
$ws.Cells.Item(18, 4).Value = 38
$ws.Cells.Item(18, 5).Value = 'Alejandro M S, Honorio G, Rigoberto Al-B, Rogelio M'

$ws.Cells.Item(19, 1).Value = 45910
$ws.Cells.Item(19, 2).Value = 'GU Henle'
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = 28.5
$ws.Cells.Item(19, 5).Value = 'Fernando V, Laurentino, Oscar VS'

$ws.Cells.Item(20, 1).Value = 45910
$ws.Cells.Item(20, 2).Value = 'HanoverSpring'
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 4).Value = 48
$ws.Cells.Item(20, 5).Value = 'Carlos Al-V, Gaudencio B, Jose P, Moises P, Rata F'

$ws.Cells.Item(21, 1).Value = 45910
$ws.Cells.Item(21, 2).Value = 'Rowan'
$ws.Cells.Item(21, 3).Value = 7
$ws.Cells.Item(21, 4).Value = 66.5
$ws.Cells.Item(21, 5).Value = 'Adalberto T, Daniel LG, Elvis T, Julio M, Luis Enrique R, Luis Martin R, Trinidad T'

$ws.Cells.Item(22, 1).Value = 45910
$ws.Cells.Item(22, 2).Value = 'Tidal Basin'
$ws.Cells.Item(22, 3).Value = 6
$ws.Cells.Item(22, 4).Value = 48
$ws.Cells.Item(22, 5).Value = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR'

$ws.Cells.Item(23, 1).Value = 45910
$ws.Cells.Item(23, 2).Value = 'Wardman'
$ws.Cells.Item(23, 3).Value = 12
$ws.Cells.Item(23, 4).Value = 114
$ws.Cells.Item(23, 5).Value = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Juan G, Miguel A, Pablo G, William A'

$ws.Cells.Item(24, 1).Value = 45911
$ws.Cells.Item(24, 2).Value = '2011 Crystal'
$ws.Cells.Item(24, 3).Value = 4
$ws.Cells.Item(24, 4).Value = 36
$ws.Cells.Item(24, 5).Value = 'Alejandro M S, Gerardo D, Rigoberto Al-B, Rogelio M'

$ws.Cells.Item(25, 1).Value = 45911
$ws.Cells.Item(25, 2).Value = 'BridgeDist'
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = 19
$ws.Cells.Item(25, 5).Value = 'Elvis T, Evaristo A'

$ws.Cells.Item(26, 1).Value = 45911
$ws.Cells.Item(26, 2).Value = 'GU Henle'
$ws.Cells.Item(26, 3).Value = 5
$ws.Cells.Item(26, 4).Value = 47.5
$ws.Cells.Item(26, 5).Value = 'Fernando V, Jose Carlos G, Laurentino, Noe VL, Oscar VS'

$ws.Cells.Item(27, 1).Value = 45911
$ws.Cells.Item(27, 2).Value = 'HanoverSpring'
$ws.Cells.Item(27, 3).Value = 5
$ws.Cells.Item(27, 4).Value = 47.5
$ws.Cells.Item(27, 5).Value = 'Antoine F, Carlos Al-V, Jose P, Moises P, Rata F'

$ws.Cells.Item(28, 1).Value = 45911
$ws.Cells.Item(28, 2).Value = 'Kingstowne'
$ws.Cells.Item(28, 3).Value = 5
$ws.Cells.Item(28, 4).Value = 47.5
$ws.Cells.Item(28, 5).Value = 'Adalberto T, Gaudencio B, Luis Enrique R, Luis Martin R, Trinidad T'

$ws.Cells.Item(29, 1).Value = 45911
$ws.Cells.Item(29, 2).Value = 'Moorefield'
$ws.Cells.Item(29, 3).Value = 7
$ws.Cells.Item(29, 4).Value = 66.5
$ws.Cells.Item(29, 5).Value = 'Benito A, Daniel GS, Daniel LG, Eric M R, Julio M, Leobardo RL, Omar S'

$ws.Cells.Item(30, 1).Value = 45911
$ws.Cells.Item(30, 2).Value = 'Tidal Basin'
$ws.Cells.Item(30, 3).Value = 6
$ws.Cells.Item(30, 4).Value = 57
$ws.Cells.Item(30, 5).Value = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR'

$ws.Cells.Item(31, 1).Value = 45911
$ws.Cells.Item(31, 2).Value = 'Wardman'
$ws.Cells.Item(31, 3).Value = 14
$ws.Cells.Item(31, 4).Value = 130
$ws.Cells.Item(31, 5).Value = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Jose Luis H, Juan G, Miguel A, Misael M, Pablo G, William A'

$ws.Cells.Item(32, 1).Value = 45912
$ws.Cells.Item(32, 2).Value = '2011 Crystal'
$ws.Cells.Item(32, 3).Value = 3
$ws.Cells.Item(32, 4).Value = 27
$ws.Cells.Item(32, 5).Value = 'Alejandro M S, Gerardo D, Rigoberto Al-B'

$ws.Cells.Item(33, 1).Value = 45912
$ws.Cells.Item(33, 2).Value = '2011 Crystal    Yard'
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 10
$ws.Cells.Item(33, 5).Value = 'Rogelio M'

$ws.Cells.Item(34, 1).Value = 45912
$ws.Cells.Item(34, 2).Value = 'Canvas, Yard'
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 9.5
$ws.Cells.Item(34, 5).Value = 'Honorio G'

$ws.Cells.Item(35, 1).Value = 45912
$ws.Cells.Item(35, 2).Value = 'GU Henle'
$ws.Cells.Item(35, 3).Value = 5
$ws.Cells.Item(35, 4).Value = 47.5
$ws.Cells.Item(35, 5).Value = 'Fernando V, Jose Carlos G, Laurentino, Noe VL, Oscar VS'

$ws.Cells.Item(36, 1).Value = 45912
$ws.Cells.Item(36, 2).Value = 'HanoverSpring'
$ws.Cells.Item(36, 3).Value = 5
$ws.Cells.Item(36, 4).Value = 47.5
$ws.Cells.Item(36, 5).Value = 'Antoine F, Carlos Al-V, Jose P, Moises P, Rata F'

$ws.Cells.Item(37, 1).Value = 45912
$ws.Cells.Item(37, 2).Value = 'Kingstowne'
$ws.Cells.Item(37, 3).Value = 8
$ws.Cells.Item(37, 4).Value = 78.5
$ws.Cells.Item(37, 5).Value = 'Adalberto T, Elvis T, Gaudencio B, Julio M, Leobardo RL, Luis Enrique R, Luis Martin R, Trinidad T'

$ws.Cells.Item(38, 1).Value = 45912
$ws.Cells.Item(38, 2).Value = 'Moorefield'
$ws.Cells.Item(38, 3).Value = 7
$ws.Cells.Item(38, 4).Value = 66.5
$ws.Cells.Item(38, 5).Value = 'Benito A, Daniel GS, Daniel LG, Eric M R, Evaristo A, Henry G, Omar S'

$ws.Cells.Item(39, 1).Value = 45912
$ws.Cells.Item(39, 2).Value = 'Tidal Basin'
$ws.Cells.Item(39, 3).Value = 7
$ws.Cells.Item(39, 4).Value = 66.5
$ws.Cells.Item(39, 5).Value = 'Alberto R, Danis BA, Eduardo H, Feliciano R, Isidro M, Juan HR, William A'

$ws.Cells.Item(40, 1).Value = 45912
$ws.Cells.Item(40, 2).Value = 'Wardman'
$ws.Cells.Item(40, 3).Value = 12
$ws.Cells.Item(40, 4).Value = 110.5
$ws.Cells.Item(40, 5).Value = 'Alfonso D, Andres G, Benny S, Carlos G, Cristobal L, Diego R, Eliacim R, Jesus L, Juan G, Miguel A, Misael M, Pablo G'

